$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 999.6667
$ws.Range("I32").Value = 1000
$ws.Range("J32").Value = 999.5
$ws.Range("K32").Value = 1000
$ws.Range("L32").Value = 999.5
$ws.Range("M32").Value = -674
$ws.Range("N32").Value = -1651.5

$ws.Range("H43").Value = 1116.7273
$ws.Range("I43").Value = 1095
$ws.Range("K43").Value = 1095
$ws.Range("M43").Value = -1026

$ws.Range("H93").Value = 27000
$ws.Range("J93").Value = 27000
$ws.Range("L93").Value = 27000
$ws.Range("N93").Value = -31992

$ws.Range("H106").Value = 7754490.5
$ws.Range("I106").Value = 14494512
$ws.Range("K106").Value = 14494512
$ws.Range("M106").Value = -14493881

$ws.Range("H129").Value = 334619.34
$ws.Range("I129").Value = 466.66666
$ws.Range("K129").Value = 1399.99998
$ws.Range("M129").Value = 3600.00002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6912
$ws.Range("I32").Value = 4591.5186
$ws.Range("K32").Value = 4591.5186
$ws.Range("M32").Value = -4304.5186

$ws.Range("H45").Value = 2260.2222
$ws.Range("I45").Value = 1790.4445
$ws.Range("J45").Value = 2730
$ws.Range("K45").Value = 1790.4445
$ws.Range("L45").Value = 2730
$ws.Range("M45").Value = -1413.4445
$ws.Range("N45").Value = -3484

$ws.Range("H61").Value = 1492.9
$ws.Range("I61").Value = 1441.6471
$ws.Range("K61").Value = 1441.6471
$ws.Range("M61").Value = -1229.6471

$ws.Range("H122").Value = 2188.125
$ws.Range("I122").Value = 1667.4
$ws.Range("J122").Value = 9999
$ws.Range("K122").Value = 5002.200000000001
$ws.Range("L122").Value = 29997
$ws.Range("M122").Value = -2552.200000000001
$ws.Range("N122").Value = -34897

$ws.Range("H132").Value = 16183.972
$ws.Range("I132").Value = 1681.963
$ws.Range("J132").Value = 65128.25
$ws.Range("K132").Value = 5045.889
$ws.Range("L132").Value = 195384.75
$ws.Range("M132").Value = -2515.889
$ws.Range("N132").Value = -200444.75

$ws.Range("H136").Value = 1492.9
$ws.Range("I136").Value = 1441.6471
$ws.Range("K136").Value = 4324.9413
$ws.Range("M136").Value = -1774.9413

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 29000
$ws.Range("J92").Value = 29000
$ws.Range("L92").Value = 29000
$ws.Range("N92").Value = -33992

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

$ws.Range("H86").Value = 27566.572
$ws.Range("I86").Value = 12475
$ws.Range("K86").Value = 12475
$ws.Range("M86").Value = -11352

$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()

$ws.Range("H89").Value = 27566.572
$ws.Range("I89").Value = 12475
$ws.Range("K89").Value = 62375
$ws.Range("M89").Value = -56759

$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()

$ws.Range("H94").Value = 3067.5
$ws.Range("I94").Value = 625
$ws.Range("J94").Value = 3881.6667
$ws.Range("K94").Value = 625
$ws.Range("L94").Value = 3881.6667
$ws.Range("M94").Value = -174
$ws.Range("N94").Value = -4783.6667

$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()

$ws.Range("H99").Value = 3475.0435
$ws.Range("I99").Value = 2992.6667
$ws.Range("K99").Value = 2992.6667
$ws.Range("M99").Value = -1494.6667

$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

$ws.Range("H103").Value = 0
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("M103").ClearContents()
$ws.Range("N103").ClearContents()

$ws.Range("H126").Value = 3475.0435
$ws.Range("I126").Value = 2992.6667
$ws.Range("K126").Value = 8978.000100000001
$ws.Range("M126").Value = -6508.000100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 720.23
$ws.Range("J131").Value = 724.4184
$ws.Range("L131").Value = 2173.2552
$ws.Range("N131").Value = -12253.2552

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 5000
$ws.Range("I43").Value = 5000
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 5000
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -4849
$ws.Range("N43").ClearContents()

$ws.Range("H46").Value = 20000
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

$ws.Range("H57").Value = 26633.334
$ws.Range("J57").Value = 26633.334
$ws.Range("L57").Value = 26633.334
$ws.Range("N57").Value = -28273.334

$ws.Range("H80").Value = 4713
$ws.Range("J80").Value = 5207
$ws.Range("L80").Value = 5207
$ws.Range("N80").Value = -7203

$ws.Range("H83").Value = 4713
$ws.Range("J83").Value = 5207
$ws.Range("L83").Value = 26035
$ws.Range("N83").Value = -36019

$ws.Range("H102").Value = 29415160
$ws.Range("J102").Value = 3566
$ws.Range("L102").Value = 3566
$ws.Range("N102").Value = -6810

$ws.Range("H132").Value = 27474.6
$ws.Range("J132").Value = 58493.223
$ws.Range("L132").Value = 175479.669
$ws.Range("N132").Value = -180539.669

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4533.067
$ws.Range("I7").Value = 2752.0952
$ws.Range("J7").Value = 8688.666999999999
$ws.Range("K7").Value = 2752.0952
$ws.Range("L7").Value = 8688.666999999999
$ws.Range("M7").Value = -2640.0952
$ws.Range("N7").Value = -8912.666999999999

$ws.Range("H94").Value = 22299.5
$ws.Range("J94").Value = 22299
$ws.Range("L94").Value = 22299
$ws.Range("N94").Value = -23651

$ws.Range("H126").Value = 4533.067
$ws.Range("I126").Value = 2752.0952
$ws.Range("J126").Value = 8688.666999999999
$ws.Range("K126").Value = 8256.285600000001
$ws.Range("L126").Value = 26066.001
$ws.Range("M126").Value = -5786.285600000001
$ws.Range("N126").Value = -31006.001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 15010000
$ws.Range("J70").Value = 15010000
$ws.Range("L70").Value = 15010000
$ws.Range("N70").Value = -15010630

$ws.Range("H73").Value = 15010000
$ws.Range("J73").Value = 15010000
$ws.Range("L73").Value = 15010000
$ws.Range("N73").Value = -15012184

$ws.Range("H107").Value = 1365.4445
$ws.Range("I107").Value = 797.8
$ws.Range("J107").Value = 2075
$ws.Range("K107").Value = 2393.4
$ws.Range("L107").Value = 6225
$ws.Range("M107").Value = -473.3999999999996
$ws.Range("N107").Value = -10065

$ws.Range("H132").Value = 1043.7073
$ws.Range("I132").Value = 614.9697
$ws.Range("K132").Value = 1844.9091
$ws.Range("M132").Value = 685.0909000000001

$ws.Range("H136").Value = 20835520
$ws.Range("I136").Value = 29412874
$ws.Range("K136").Value = 88238622
$ws.Range("M136").Value = -88236072
